$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Unmerge the affected area so per-cell writes land on the correct cells
#    (writing into a non-anchor cell of a merged range is silently ignored).
$ws.Range("A1:N40").UnMerge()

# 2) Shift rows 22..33 down to 23..34 (bottom-up so we don't clobber source
#    data before it's copied), making room for the new item row at 22.
#    Column A (the running item-number counter) is left untouched since it
#    already holds the correct sequential value for every row number.
for ($r = 33; $r -ge 22; $r--) {
    $src = $ws.Range("B" + $r + ":N" + $r)
    $dst = $ws.Range("B" + ($r + 1) + ":N" + ($r + 1))
    $dst.Value2 = $src.Value2
}
$ws.Range("A32").Value2 = 29

# 3) Fill the new row 22 with the PENDULINE item's data.
$ws.Range("A22").Value2 = 19
$ws.Range("B22").Value2 = "PENDULINE كريم ب زبدة الشيا"
$ws.Range("C22").Value2 = ""
$ws.Range("D22").Value2 = ""
$ws.Range("E22").Value2 = ""
$ws.Range("F22").Value2 = ""
$ws.Range("G22").Value2 = ""
$ws.Range("H22").Value2 = "1:0"
$ws.Range("I22").Value2 = ""
$ws.Range("J22").Value2 = ""
$ws.Range("K22").Value2 = ""
$ws.Range("L22").Value2 = 160.05
$ws.Range("M22").Value2 = ""
$ws.Range("N22").Value2 = "1:0"

# 4) Update the grand-total cell (was row 32, now row 33) to include the new row.
$ws.Range("K33").Value2 = 2331.63

# 5) Re-create the merged cells: per-item rows 4..32, the totals row (33) and
#    the footer row (34).
for ($r = 4; $r -le 32; $r++) {
    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
}
$ws.Range("K33:N33").Merge()
$ws.Range("A34:E34").Merge()
$ws.Range("F34:G34").Merge()
$ws.Range("I34:N34").Merge()

# 6) Row heights: the inserted row and the rows that used to carry the
#    totals/footer content pick up slightly different heights once Excel
#    reflows them.
$ws.Rows.Item(32).RowHeight = 25.5
$ws.Rows.Item(33).RowHeight = 25.5
$ws.Rows.Item(34).RowHeight = 17.25
